$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text NumberFormat to Price (D) cells whose new values would
# otherwise be auto-parsed as numbers by Excel, so they remain text
# like the original inline-string cells.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

# Apply the updated cell values from the diff.

# Row 2
$ws.Range("D2").Value = "69.160.50"
$ws.Range("E2").Value = "  +2.36%  "

# Row 3
$ws.Range("D3").Value = "3.777.94"
$ws.Range("E3").Value = "  +0.32%  "

# Row 4
$ws.Range("E4").Value = "  -0.49%  "

# Row 5
$ws.Range("D5").Value = "625.19"
$ws.Range("E5").Value = "  +4.46%  "

# Row 6
$ws.Range("D6").Value = "166.23"
$ws.Range("E6").Value = "  +2.31%  "

# Row 7
$ws.Range("D7").Value = "3.776.44"
$ws.Range("E7").Value = "  +0.37%  "

# Row 8
$ws.Range("E8").Value = "  -0.11%  "

# Row 9
$ws.Range("E9").Value = "  +1.78%  "

# Row 10
$ws.Range("E10").Value = "  +3.32%  "

# Row 11
$ws.Range("E11").Value = "  +2.96%  "

# Row 12
$ws.Range("E12").Value = "  +2.07%  "

# Row 13
$ws.Range("D13").Value = "0.0000247"
$ws.Range("E13").Value = "  +1.40%  "

# Row 14
$ws.Range("D14").Value = "35.64"
$ws.Range("E14").Value = "  +1.94%  "

# Row 15
$ws.Range("D15").Value = "4.413.35"
$ws.Range("E15").Value = "  +0.34%  "

# Row 16
$ws.Range("D16").Value = "3.789.57"
$ws.Range("E16").Value = "  -0.09%  "

# Row 17
$ws.Range("D17").Value = "69.162.27"
$ws.Range("E17").Value = "  +2.33%  "

# Row 18
$ws.Range("D18").Value = "17.66"
$ws.Range("E18").Value = "  -2.66%  "

# Row 19
$ws.Range("E19").Value = "  +2.02%  "

# Row 20
$ws.Range("E20").Value = "  -0.89%  "

# Row 21
$ws.Range("D21").Value = "468.40"
$ws.Range("E21").Value = "  +2.60%  "

# Row 22
$ws.Range("D22").Value = "9.63"
$ws.Range("E22").Value = "  +2.20%  "

# Row 23
$ws.Range("D23").Value = "0.707"
$ws.Range("E23").Value = "  +2.65%  "

# Row 24
$ws.Range("D24").Value = "0.0000148"
$ws.Range("E24").Value = "  +4.53%  "

# Row 25
$ws.Range("D25").Value = "83.15"
$ws.Range("E25").Value = "  +0.71%  "

# Row 26
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D26").Value = "12.05"
$ws.Range("E26").Value = "  +2.01%  "

# Row 27
$ws.Range("B27").Value = "Fetch.AI"
$ws.Range("C27").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D27").Value = "2.16"
$ws.Range("E27").Value = "  +4.21%  "

# Row 28
$ws.Range("D28").Value = "10.01"
$ws.Range("E28").Value = "  +2.26%  "

# Row 29
$ws.Range("E29").Value = "  -0.06%  "

# Row 30
$ws.Range("D30").Value = "3.926.18"
$ws.Range("E30").Value = "  +0.27%  "

# Row 31
$ws.Range("E31").Value = "  +3.89%  "

# Row 32
$ws.Range("E32").Value = "  +2.25%  "

# Row 33
$ws.Range("D33").Value = "7.24"
$ws.Range("E33").Value = "  +1.12%  "

# Row 34
$ws.Range("D34").Value = "28.77"
$ws.Range("E34").Value = "  +0.07%  "

# Row 35
$ws.Range("E35").Value = "  -0.21%  "

# Row 36
$ws.Range("E36").Value = "  +16.25%  "

# Row 37
$ws.Range("E37").Value = "  +1.05%  "

# Row 38
$ws.Range("D38").Value = "3.728.23"
$ws.Range("E38").Value = "  +0.34%  "

# Row 39
$ws.Range("E39").Value = "  +3.05%  "

# Row 40
$ws.Range("E40").Value = "  +8.79%  "

# Row 41
$ws.Range("D41").Value = "5.81"
$ws.Range("E41").Value = "  +1.25%  "

# Row 42
$ws.Range("D42").Value = "0.969"
$ws.Range("E42").Value = "  -0.66%  "

# Row 43
$ws.Range("E43").Value = "  +0.05%  "

# Row 44
$ws.Range("E44").Value = "  -0.01%  "

# Row 45
$ws.Range("E45").Value = "  +1.54%  "

# Row 46
$ws.Range("D46").Value = "43.01"
$ws.Range("E46").Value = "  -0.42%  "

# Row 47
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "152.55"
$ws.Range("E47").Value = "  +0.42%  "

# Row 48
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "1.92"
$ws.Range("E48").Value = "  +5.21%  "

# Row 49
$ws.Range("E49").Value = "  -0.85%  "

# Row 50
$ws.Range("E50").Value = "  +1.97%  "

# Row 51
$ws.Range("E51").Value = "  -0.07%  "
